# Update NATMI LR-pair stats (Adam17-Itgb1) for rows 2-10 following Dr Hou's advice:
# ligand/receptor expressing-cell counts changed 1 -> 3, with corresponding
# recalculated expression/specificity values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.47808166666666
$ws.Range("H2").Value = 100.434245
$ws.Range("I2").Value = 0.4880542983452505
$ws.Range("J2").Value = 0.4880542983452505
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 3766.732525969679
$ws.Range("R2").Value = 33900.59273372711
$ws.Range("S2").Value = 0.1598530577030708
$ws.Range("T2").Value = 0.1598530577030708

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.47808166666666
$ws.Range("H3").Value = 100.434245
$ws.Range("I3").Value = 0.4880542983452505
$ws.Range("J3").Value = 0.4880542983452505
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 3559.204375096056
$ws.Range("R3").Value = 32032.83937586451
$ws.Range("S3").Value = 0.1510459525402021
$ws.Range("T3").Value = 0.1510459525402021

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.47808166666666
$ws.Range("H4").Value = 100.434245
$ws.Range("I4").Value = 0.4880542983452505
$ws.Range("J4").Value = 0.4880542983452505
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 4174.437420401186
$ws.Range("R4").Value = 37569.93678361067
$ws.Range("S4").Value = 0.1771552881019776
$ws.Range("T4").Value = 0.1771552881019776

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.24776266666667
$ws.Range("H5").Value = 69.743288
$ws.Range("I5").Value = 0.3389134003957588
$ws.Range("J5").Value = 0.3389134003957588
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 2615.684634037632
$ws.Range("R5").Value = 23541.16170633869
$ws.Range("S5").Value = 0.1110047458520337
$ws.Range("T5").Value = 0.1110047458520337

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.24776266666667
$ws.Range("H6").Value = 69.743288
$ws.Range("I6").Value = 0.3389134003957588
$ws.Range("J6").Value = 0.3389134003957588
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 2471.573473601403
$ws.Range("R6").Value = 22244.16126241262
$ws.Range("S6").Value = 0.1048889387205096
$ws.Range("T6").Value = 0.1048889387205096

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.24776266666667
$ws.Range("H7").Value = 69.743288
$ws.Range("I7").Value = 0.3389134003957588
$ws.Range("J7").Value = 0.3389134003957588
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 2898.802009703135
$ws.Range("R7").Value = 26089.21808732821
$ws.Range("S7").Value = 0.1230197158232155
$ws.Range("T7").Value = 0.1230197158232155

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.86914966666667
$ws.Range("H8").Value = 35.607449
$ws.Range("I8").Value = 0.1730323012589908
$ws.Range("J8").Value = 0.1730323012589908
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 1335.438289152336
$ws.Range("R8").Value = 12018.94460237102
$ws.Range("S8").Value = 0.05667349418175195
$ws.Range("T8").Value = 0.05667349418175194

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.86914966666667
$ws.Range("H9").Value = 35.607449
$ws.Range("I9").Value = 0.1730323012589908
$ws.Range("J9").Value = 0.1730323012589908
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 1261.862308685744
$ws.Range("R9").Value = 11356.7607781717
$ws.Range("S9").Value = 0.05355106768345465
$ws.Range("T9").Value = 0.05355106768345465

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.86914966666667
$ws.Range("H10").Value = 35.607449
$ws.Range("I10").Value = 0.1730323012589908
$ws.Range("J10").Value = 0.1730323012589908
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 1479.983919335748
$ws.Range("R10").Value = 13319.85527402173
$ws.Range("S10").Value = 0.06280773939378421
$ws.Range("T10").Value = 0.06280773939378421
